$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ENGLISH"
$ws.Range("R2").ClearContents()
